$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 167; existing rows 167:184 shift down to 168:185
$ws.Rows.Item(167).Insert()

# Populate the newly inserted row 167 with the new weekly record
$ws.Range("A167").Value = 8
$ws.Range("B167").Value = "Terminal La Palmera de La Serena"
$ws.Range("C167").Value = "Coquimbo"
$ws.Range("D167").Value = 44449
$ws.Range("E167").Value = 4
$ws.Range("F167").Value = 100114013
$ws.Range("G167").Value = "Zanahoria"
$ws.Range("H167").Value = "Sin especificar"
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 700
$ws.Range("K167").Value = 4500
$ws.Range("L167").Value = 5000
$ws.Range("M167").Value = 4750
$ws.Range("N167").Value = "`$/saco 20 kilos"
$ws.Range("O167").Value = "Provincia del Elquí"
$ws.Range("P167").Value = 238
$ws.Range("Q167").Value = 20
$ws.Range("R167").Value = "Hortaliza"
